# Applies scheduled-runner value updates to the Leve profit sheets.
# Generated from the authoritative cell-level diff (ALC, ARM, BSM, CRP, CUL, GSM).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 5556439
$ws.Range("J17").Value = 6250931.5
$ws.Range("L17").Value = 18752794.5
$ws.Range("N17").Value = -18753130.5
$ws.Range("H19").Value = 809.9474
$ws.Range("I19").Value = 890
$ws.Range("J19").Value = 699.875
$ws.Range("K19").Value = 890
$ws.Range("L19").Value = 699.875
$ws.Range("M19").Value = -715
$ws.Range("N19").Value = -1049.875
$ws.Range("H32").Value = 1326.5
$ws.Range("I32").Value = 700
$ws.Range("J32").Value = 1430.9166
$ws.Range("K32").Value = 700
$ws.Range("L32").Value = 1430.9166
$ws.Range("M32").Value = -374
$ws.Range("N32").Value = -2082.9166
$ws.Range("H33").Value = 4156.28
$ws.Range("I33").Value = 5650.778
$ws.Range("K33").Value = 5650.778
$ws.Range("M33").Value = -5421.778
$ws.Range("H40").Value = 1829.9459
$ws.Range("I40").Value = 1596.1765
$ws.Range("J40").Value = 2028.65
$ws.Range("K40").Value = 1596.1765
$ws.Range("L40").Value = 2028.65
$ws.Range("M40").Value = -1421.1765
$ws.Range("N40").Value = -2378.65
$ws.Range("H51").Value = 4673
$ws.Range("I51").Value = 2542.25
$ws.Range("K51").Value = 2542.25
$ws.Range("M51").Value = -2058.25
$ws.Range("H86").Value = 102100
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 102100
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 102100
$ws.Range("M86").ClearContents()  # was -1477
$ws.Range("N86").Value = -104346
$ws.Range("H89").Value = 102100
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 102100
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 510500
$ws.Range("M89").ClearContents()  # was -7384
$ws.Range("N89").Value = -521732
$ws.Range("H112").Value = 1320.6552
$ws.Range("I112").Value = 1280
$ws.Range("J112").Value = 1325.3462
$ws.Range("K112").Value = 3840
$ws.Range("L112").Value = 3976.0386
$ws.Range("M112").Value = -2732
$ws.Range("N112").Value = -6192.0386
$ws.Range("H113").Value = 8066.087
$ws.Range("I113").Value = 2347.6924
$ws.Range("K113").Value = 2347.6924
$ws.Range("M113").Value = 906.3076000000001
$ws.Range("H138").Value = 1770.8518
$ws.Range("I138").Value = 1426.2188
$ws.Range("K138").Value = 4278.6564
$ws.Range("M138").Value = 861.3436000000002

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3018.75
$ws.Range("I45").Value = 3134.6667
$ws.Range("J45").Value = 2671
$ws.Range("K45").Value = 3134.6667
$ws.Range("L45").Value = 2671
$ws.Range("M45").Value = -2757.6667
$ws.Range("N45").Value = -3425
$ws.Range("H97").Value = 896.4286
$ws.Range("I97").Value = 916.913
$ws.Range("J97").Value = 802.2
$ws.Range("K97").Value = 916.913
$ws.Range("L97").Value = 802.2
$ws.Range("M97").Value = -420.913
$ws.Range("N97").Value = -1794.2
$ws.Range("H102").Value = 2047.8572
$ws.Range("I102").Value = 2030.8334
$ws.Range("K102").Value = 2030.8334
$ws.Range("M102").Value = -408.8334

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 556.02563
$ws.Range("I94").Value = 403.68182
$ws.Range("J94").Value = 753.17645
$ws.Range("K94").Value = 403.68182
$ws.Range("L94").Value = 753.17645
$ws.Range("M94").Value = 47.31817999999998
$ws.Range("N94").Value = -1655.17645
$ws.Range("H105").Value = 2838.75
$ws.Range("I105").Value = 2838.75
$ws.Range("K105").Value = 2838.75
$ws.Range("M105").Value = -1091.75
$ws.Range("H107").Value = 551.6429000000001
$ws.Range("I107").Value = 390
$ws.Range("K107").Value = 390
$ws.Range("M107").Value = 1530

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 750.34784
$ws.Range("I16").Value = 748.0909
$ws.Range("J16").Value = 800
$ws.Range("K16").Value = 748.0909
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = -461.0909
$ws.Range("N16").Value = -1374
$ws.Range("H86").Value = 4855.875
$ws.Range("I86").Value = 4307.8335
$ws.Range("K86").Value = 4307.8335
$ws.Range("M86").Value = -3184.8335
$ws.Range("H89").Value = 4855.875
$ws.Range("I89").Value = 4307.8335
$ws.Range("K89").Value = 21539.1675
$ws.Range("M89").Value = -15923.1675
$ws.Range("H105").Value = 1795.0344
$ws.Range("I105").Value = 1991
$ws.Range("K105").Value = 1991
$ws.Range("M105").Value = -244
$ws.Range("H113").Value = 750.34784
$ws.Range("I113").Value = 748.0909
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 748.0909
$ws.Range("L113").Value = 800
$ws.Range("M113").Value = 1421.9091
$ws.Range("N113").Value = -5140

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 46
$ws.Range("I12").Value = 5
$ws.Range("K12").Value = 15
$ws.Range("M12").Value = 158
$ws.Range("H23").Value = 84.25
$ws.Range("J23").Value = 84.25
$ws.Range("L23").Value = 252.75
$ws.Range("N23").Value = -722.75
$ws.Range("H69").Value = 2970.6667
$ws.Range("I69").Value = 5012
$ws.Range("J69").Value = 1950
$ws.Range("K69").Value = 15036
$ws.Range("L69").Value = 5850
$ws.Range("M69").Value = -14225
$ws.Range("N69").Value = -7472
$ws.Range("H72").Value = 2970.6667
$ws.Range("I72").Value = 5012
$ws.Range("J72").Value = 1950
$ws.Range("K72").Value = 45108
$ws.Range("L72").Value = 17550
$ws.Range("M72").Value = -41052
$ws.Range("N72").Value = -25662
$ws.Range("H74").Value = 5750
$ws.Range("J74").Value = 5750
$ws.Range("L74").Value = 17250
$ws.Range("N74").Value = -19372
$ws.Range("H77").Value = 5750
$ws.Range("J77").Value = 5750
$ws.Range("L77").Value = 51750
$ws.Range("N77").Value = -62358
$ws.Range("H80").Value = 4381.7
$ws.Range("I80").Value = 1999.5
$ws.Range("J80").Value = 4977.25
$ws.Range("K80").Value = 5998.5
$ws.Range("L80").Value = 14931.75
$ws.Range("M80").Value = -5062.5
$ws.Range("N80").Value = -16803.75
$ws.Range("H83").Value = 4381.7
$ws.Range("I83").Value = 1999.5
$ws.Range("J83").Value = 4977.25
$ws.Range("K83").Value = 17995.5
$ws.Range("L83").Value = 44795.25
$ws.Range("M83").Value = -13315.5
$ws.Range("N83").Value = -54155.25
$ws.Range("H92").Value = 1073.909
$ws.Range("I92").Value = 1183
$ws.Range("K92").Value = 3549
$ws.Range("M92").Value = -2301
$ws.Range("H131").Value = 1414.4025
$ws.Range("J131").Value = 1532.6198
$ws.Range("L131").Value = 4597.859399999999
$ws.Range("N131").Value = -14677.8594

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2026.3158
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 2083.3333
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2083.3333
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -4079.3333
$ws.Range("H83").Value = 2026.3158
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 2083.3333
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 10416.6665
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -20400.6665
